$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Updated cryptos list - refresh of price (D) and volume/1h change (E)
# columns, plus a few rows that got re-ordered (B/C/D/E all changed).
#
# Several "Price" values are strings that look like plain numbers
# (e.g. "525.94"); Excel's Range.Value setter auto-converts those to
# numeric cells, but the source file stores them as plain text. Force
# text storage for those by temporarily marking the cell as Text
# before assigning, then restoring the original (default) style so the
# on-disk style index is unaffected.
# ----------------------------------------------------------------------

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.817.33"
$ws.Range("E2").Value = "  +1.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.052.34"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "525.94"
$ws.Range("E5").Value = "  +5.46%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "142.44"
$ws.Range("E6").Value = "  +4.93%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.449"
$ws.Range("E8").Value = "  +5.09%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +5.51%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +7.43%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +5.61%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.34%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.577.12"
$ws.Range("E13").Value = "  +2.08%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "27.08"
$ws.Range("E14").Value = "  +7.72%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +15.59%  "

# Row 16 & 17 - WrappedBTC and Polkadot swap order
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D16") "6.30"
$ws.Range("E16").Value = "  +7.85%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "57.773.90"
$ws.Range("E17").Value = "  +2.08%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.051.93"
$ws.Range("E18").Value = "  +2.18%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "13.10"
$ws.Range("E19").Value = "  +5.31%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "8.17"
$ws.Range("E20").Value = "  +4.98%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "338.10"
$ws.Range("E21").Value = "  +3.50%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.09%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +7.02%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "65.00"
$ws.Range("E24").Value = "  +5.53%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  +6.16%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "0.0₃0973"
$ws.Range("E26").Value = "  +7.58%  "

# Row 27 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D27") "0.998"
$ws.Range("E27").Value = "  +0.05%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  +6.04%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "7.36"
$ws.Range("E29").Value = "  +9.27%  "

# Row 30 - PancakeSwap
Set-TextValue $ws.Range("D30") "1.85"
$ws.Range("E30").Value = "  +6.42%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +4.47%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "21.13"
$ws.Range("E32").Value = "  +4.92%  "

# Row 33 - NEARProtocol
Set-TextValue $ws.Range("D33") "4.74"
$ws.Range("E33").Value = "  +5.46%  "

# Row 34 - Monero
Set-TextValue $ws.Range("D34") "156.15"
$ws.Range("E34").Value = "  +1.52%  "

# Row 35 - Aptos
Set-TextValue $ws.Range("D35") "5.99"
$ws.Range("E35").Value = "  +6.75%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +2.79%  "

# Row 37 - EnergySwap
Set-TextValue $ws.Range("D37") "26.16"
$ws.Range("E37").Value = "  +12.64%  "

# Row 38 - Hedera
Set-TextValue $ws.Range("D38") "0.0702"
$ws.Range("E38").Value = "  +3.02%  "

# Row 39 - RenzoRestakedETH
$ws.Range("D39").Value = "3.089.13"
$ws.Range("E39").Value = "  +2.03%  "

# Row 40 - OKB
Set-TextValue $ws.Range("D40") "37.74"
$ws.Range("E40").Value = "  +2.98%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +8.32%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.24%  "

# Row 43 - Stacks
Set-TextValue $ws.Range("D43") "1.47"
$ws.Range("E43").Value = "  +4.65%  "

# Row 44 - Mantle
Set-TextValue $ws.Range("D44") "0.662"
$ws.Range("E44").Value = "  +3.30%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.326.68"
$ws.Range("E45").Value = "  +4.00%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  +3.03%  "

# Row 47, 48 & 49 - dogwifhat, VeChain, Cosmos re-ordered
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0246"
$ws.Range("E47").Value = "  +3.13%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "6.05"
$ws.Range("E48").Value = "  +4.70%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D49") "1.99"
$ws.Range("E49").Value = "  +2.05%  "

# Row 50 - InjectiveProtocol
Set-TextValue $ws.Range("D50") "20.08"
$ws.Range("E50").Value = "  +4.83%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +6.40%  "
